$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a plain text value to a cell, forcing text storage so that
# numeric-looking strings (e.g. "1.0000", "316.42") are NOT silently
# converted into real numbers by Excel, and without leaving a stray
# number-format style applied to the cell afterwards.
function Set-TextValue($sheet, $addr, $val) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "28.063.01"
Set-TextValue $ws "E2" "  -1.16%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "1.791.49"
Set-TextValue $ws "E3" "  -0.72%  "

# Row 4 - TetherUSD
Set-TextValue $ws "D4" "1.0000"
Set-TextValue $ws "E4" "  -0.06%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "316.42"
Set-TextValue $ws "E5" "  +0.06%  "

# Row 6 - USDC
Set-TextValue $ws "D6" "0.9989"
Set-TextValue $ws "E6" "  -0.13%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.5340"
Set-TextValue $ws "E7" "  -2.61%  "

# Row 8 - Cardano
Set-TextValue $ws "D8" "0.3764"
Set-TextValue $ws "E8" "  -2.18%  "

# Row 9 - Dogecoin
Set-TextValue $ws "D9" "0.07450"
Set-TextValue $ws "E9" "  -1.74%  "

# Row 10 - OKB
Set-TextValue $ws "D10" "41.80"
Set-TextValue $ws "E10" "  -1.43%  "

# Row 11 - Polygon
Set-TextValue $ws "D11" "1.096"
Set-TextValue $ws "E11" "  -2.68%  "

# Row 12 - BinanceUSD
Set-TextValue $ws "D12" "0.9976"
Set-TextValue $ws "E12" "  -0.29%  "

# Row 13 - Solana
Set-TextValue $ws "D13" "20.64"
Set-TextValue $ws "E13" "  -2.58%  "

# Row 14 - Polkadot
Set-TextValue $ws "D14" "6.108"
Set-TextValue $ws "E14" "  -1.26%  "

# Row 15 - was Chainlink, now WrappedEther (rows 15/16 swapped)
Set-TextValue $ws "B15" "WrappedEther"
Set-TextValue $ws "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D15" "1.795.38"
Set-TextValue $ws "E15" "  -0.61%  "

# Row 16 - was WrappedEther, now Chainlink
Set-TextValue $ws "B16" "Chainlink"
Set-TextValue $ws "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws "D16" "7.232"
Set-TextValue $ws "E16" "  -2.12%  "

# Row 17 - Litecoin
Set-TextValue $ws "E17" "  -3.48%  "

# Row 18 - ShibaInu
Set-TextValue $ws "D18" "0.00001054"
Set-TextValue $ws "E18" "  -1.52%  "

# Row 19 - TRON
Set-TextValue $ws "D19" "0.06451"
Set-TextValue $ws "E19" "  +0.05%  "

# Row 20 - Dai
Set-TextValue $ws "D20" "0.9984"
Set-TextValue $ws "E20" "  -0.13%  "

# Row 21 - Avalanche
Set-TextValue $ws "D21" "17.34"
Set-TextValue $ws "E21" "  +0.04%  "

# Row 22 - Uniswap
Set-TextValue $ws "D22" "5.898"
Set-TextValue $ws "E22" "  -1.29%  "

# Row 23 - WrappedBTC
Set-TextValue $ws "D23" "28.114.79"
Set-TextValue $ws "E23" "  -1.06%  "

# Row 24 - Cosmos
Set-TextValue $ws "D24" "11.19"
Set-TextValue $ws "E24" "  -2.17%  "

# Row 25 - Toncoin
Set-TextValue $ws "D25" "2.099"
Set-TextValue $ws "E25" "  -1.49%  "

# Row 26 - Monero
Set-TextValue $ws "D26" "155.24"
Set-TextValue $ws "E26" "  -2.45%  "

# Row 27 - EthereumClassic
Set-TextValue $ws "D27" "20.25"
Set-TextValue $ws "E27" "  -1.96%  "

# Row 28 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D28" "1.984.17"
Set-TextValue $ws "E28" "  -1.48%  "

# Row 29 - LidoDAOToken
Set-TextValue $ws "D29" "2.292"
Set-TextValue $ws "E29" "  -4.57%  "

# Row 30 - BitcoinCash
Set-TextValue $ws "D30" "120.17"
Set-TextValue $ws "E30" "  -3.03%  "

# Row 31 - ImmutableX
Set-TextValue $ws "D31" "1.115"
Set-TextValue $ws "E31" "  -0.58%  "

# Row 32 - Stellar
Set-TextValue $ws "E32" "  +2.98%  "

# Row 33 - was Filecoin, now HuobiToken (rows 33/34 swapped)
Set-TextValue $ws "B33" "HuobiToken"
Set-TextValue $ws "C33" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D33" "3.642"
Set-TextValue $ws "E33" "  -1.12%  "

# Row 34 - was HuobiToken, now Filecoin
Set-TextValue $ws "B34" "Filecoin"
Set-TextValue $ws "C34" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D34" "5.553"
Set-TextValue $ws "E34" "  -3.31%  "

# Row 35 - Algorand
Set-TextValue $ws "D35" "0.2260"
Set-TextValue $ws "E35" "  -2.00%  "

# Row 36 - Hedera
Set-TextValue $ws "D36" "0.06470"
Set-TextValue $ws "E36" "  +0.34%  "

# Row 37 - VeChain
Set-TextValue $ws "E37" "  -1.61%  "

# Row 38 - InternetComputer(DFINITY)
Set-TextValue $ws "D38" "5.018"
Set-TextValue $ws "E38" "  -3.06%  "

# Row 39 - FraxShare
Set-TextValue $ws "D39" "8.524"
Set-TextValue $ws "E39" "  -3.65%  "

# Row 40 - WEMIXTOKEN
Set-TextValue $ws "E40" "  +4.50%  "

# Row 41 - TheSandbox
Set-TextValue $ws "E41" "  -4.02%  "

# Row 42 - Aptos
Set-TextValue $ws "D42" "11.07"
Set-TextValue $ws "E42" "  -4.79%  "

# Row 43 - TrustWalletToken
Set-TextValue $ws "D43" "1.172"
Set-TextValue $ws "E43" "  +0.96%  "

# Row 44 - Frax
Set-TextValue $ws "D44" "0.9977"
Set-TextValue $ws "E44" "  -0.16%  "

# Row 45 - EnergySwap
Set-TextValue $ws "D45" "13.23"
Set-TextValue $ws "E45" "  -2.00%  "

# Row 46 - PancakeSwap
Set-TextValue $ws "D46" "3.671"
Set-TextValue $ws "E46" "  -0.40%  "

# Row 47 - Decentraland
Set-TextValue $ws "D47" "0.5766"
Set-TextValue $ws "E47" "  -3.65%  "

# Row 48 - Quant
Set-TextValue $ws "D48" "127.07"
Set-TextValue $ws "E48" "  +0.11%  "

# Row 49 - EOS
Set-TextValue $ws "D49" "1.189"
Set-TextValue $ws "E49" "  +3.64%  "

# Row 50 - NEARProtocol
Set-TextValue $ws "E50" "  -2.90%  "

# Row 51 - Cronos
Set-TextValue $ws "D51" "0.06804"
Set-TextValue $ws "E51" "  -1.40%  "
